$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr1 = New-Object 'object[,]' 24,5
$arr1[0,0] = 1.02
$arr1[0,1] = 1.03003059421967
$arr1[0,2] = 1.032675008380089
$arr1[0,3] = 1.029764427772845
$arr1[0,4] = 1.039563603515629
$arr1[1,0] = 1.02
$arr1[1,1] = 1.031326642095284
$arr1[1,2] = 1.033884536087532
$arr1[1,3] = 1.030874938474305
$arr1[1,4] = 1.041017841544534
$arr1[2,0] = 1.02
$arr1[2,1] = 1.032164146635687
$arr1[2,2] = 1.034666373499147
$arr1[2,3] = 1.031592865498182
$arr1[2,4] = 1.041957721137565
$arr1[3,0] = 1.02
$arr1[3,1] = 1.03251596858578
$arr1[3,2] = 1.034994867985742
$arr1[3,3] = 1.031894529888683
$arr1[3,4] = 1.042352586218163
$arr1[4,0] = 1.02
$arr1[4,1] = 1.032575025629164
$arr1[4,2] = 1.035050012615892
$arr1[4,3] = 1.031945171830716
$arr1[4,4] = 1.042418870693522
$arr1[5,0] = 1.02
$arr1[5,1] = 1.032168848733476
$arr1[5,2] = 1.034670763600712
$arr1[5,3] = 1.031596896946696
$arr1[5,4] = 1.041962998360285
$arr1[6,0] = 1.02
$arr1[6,1] = 1.030468835277776
$arr1[6,2] = 1.033083943295358
$arr1[6,3] = 1.030139865929588
$arr1[6,4] = 1.040055303058524
$arr1[7,0] = 1.02
$arr1[7,1] = 1.027464379706436
$arr1[7,2] = 1.030281420023796
$arr1[7,3] = 1.027567298998065
$arr1[7,4] = 1.036684960641742
$arr1[8,0] = 1.02
$arr1[8,1] = 1.02545520364613
$arr1[8,2] = 1.028408585734066
$arr1[8,3] = 1.025848643487955
$arr1[8,4] = 1.03443184996885
$arr1[9,0] = 1.02
$arr1[9,1] = 1.024583673381568
$arr1[9,2] = 1.027596513906477
$arr1[9,3] = 1.02510354768312
$arr1[9,4] = 1.033454677431676
$arr1[10,0] = 1.02
$arr1[10,1] = 1.024259710776896
$arr1[10,2] = 1.0272947012557
$arr1[10,3] = 1.024826645949892
$arr1[10,4] = 1.033091470866632
$arr1[11,0] = 1.02
$arr1[11,1] = 1.024329212690006
$arr1[11,2] = 1.027359448994467
$arr1[11,3] = 1.024886048706641
$arr1[11,4] = 1.033169390929473
$arr1[12,0] = 1.02
$arr1[12,1] = 1.024556899407813
$arr1[12,2] = 1.02757156953377
$arr1[12,3] = 1.025080661776835
$arr1[12,4] = 1.033424659622514
$arr1[13,0] = 1.02
$arr1[13,1] = 1.024697153077597
$arr1[13,2] = 1.027702240939819
$arr1[13,3] = 1.025200550651068
$arr1[13,4] = 1.033581906930013
$arr1[14,0] = 1.02
$arr1[14,1] = 1.02551301124328
$arr1[14,2] = 1.028462456246803
$arr1[14,3] = 1.02589807362034
$arr1[14,4] = 1.034496668238739
$arr1[15,0] = 1.02
$arr1[15,1] = 1.026024360315759
$arr1[15,2] = 1.028939015667445
$arr1[15,3] = 1.02633536597198
$arr1[15,4] = 1.035070051274831
$arr1[16,0] = 1.02
$arr1[16,1] = 1.026322473003008
$arr1[16,2] = 1.029216876737816
$arr1[16,3] = 1.026590343982846
$arr1[16,4] = 1.035404345764915
$arr1[17,0] = 1.02
$arr1[17,1] = 1.026424096690637
$arr1[17,2] = 1.029311602000102
$arr1[17,3] = 1.026677270214442
$arr1[17,4] = 1.035518306283044
$arr1[18,0] = 1.02
$arr1[18,1] = 1.025969512808981
$arr1[18,2] = 1.028887896553569
$arr1[18,3] = 1.026288457691102
$arr1[18,4] = 1.035008548280702
$arr1[19,0] = 1.02
$arr1[19,1] = 1.024489857943978
$arr1[19,2] = 1.027509110133653
$arr1[19,3] = 1.025023356972683
$arr1[19,4] = 1.033349496066447
$arr1[20,0] = 1.02
$arr1[20,1] = 1.023558163570405
$arr1[20,2] = 1.026641209154701
$arr1[20,3] = 1.024227126194205
$arr1[20,4] = 1.032304985371216
$arr1[21,0] = 1.02
$arr1[21,1] = 1.024052204729324
$arr1[21,2] = 1.027101396462662
$arr1[21,3] = 1.02464930140389
$arr1[21,4] = 1.032858835082407
$arr1[22,0] = 1.02
$arr1[22,1] = 1.025994296495141
$arr1[22,2] = 1.02891099541413
$arr1[22,3] = 1.026309653793401
$arr1[22,4] = 1.035036339299814
$arr1[23,0] = 1.02
$arr1[23,1] = 1.028242176690631
$arr1[23,2] = 1.031006713764638
$arr1[23,3] = 1.028232991831201
$arr1[23,4] = 1.037557345448684

$ws.Range("B2:F25").Value = $arr1

$arr2 = New-Object 'object[,]' 24,6
$arr2[0,0] = 1.034293656520389
$arr2[0,1] = 1.035174091447719
$arr2[0,2] = 1.035479555945998
$arr2[0,3] = 1.032577392075639
$arr2[0,4] = 1.042348434739381
$arr2[0,5] = 1.01563566308488
$arr2[1,0] = 1.034603817565304
$arr2[1,1] = 1.036109651998183
$arr2[1,2] = 1.036496834073602
$arr2[1,3] = 1.033495303450516
$arr2[1,4] = 1.04361122057631
$arr2[1,5] = 1.015952517801043
$arr2[2,0] = 1.034802135815353
$arr2[2,1] = 1.036713456052028
$arr2[2,2] = 1.037153715486788
$arr2[2,3] = 1.034088026717336
$arr2[2,4] = 1.044426712150752
$arr2[2,5] = 1.016156833291937
$arr2[3,0] = 1.034884940984408
$arr2[3,1] = 1.036966923284954
$arr2[3,2] = 1.037429544624646
$arr2[3,3] = 1.034336916601135
$arr2[3,4] = 1.044769162787273
$arr2[3,5] = 1.016242558314284
$arr2[4,0] = 1.034898811072642
$arr2[4,1] = 1.037009459820397
$arr2[4,2] = 1.037475838652258
$arr2[4,3] = 1.034378689323229
$arr2[4,4] = 1.044826639489669
$arr2[4,5] = 1.016256942027158
$arr2[5,0] = 1.034803244491821
$arr2[5,1] = 1.036716844351633
$arr2[5,2] = 1.037157402397743
$arr2[5,3] = 1.034091353536144
$arr2[5,4] = 1.044431289488224
$arr2[5,5] = 1.016157979417611
$arr2[6,0] = 1.034398969709173
$arr2[6,1] = 1.035490594589919
$arr2[6,2] = 1.035823634478368
$arr2[6,3] = 1.032887860173135
$arr2[6,4] = 1.04277553641307
$arr2[6,5] = 1.015742893337774
$arr2[7,0] = 1.033668330458403
$arr2[7,1] = 1.033317657304211
$arr2[7,2] = 1.033462766091848
$arr2[7,3] = 1.030757632966479
$arr2[7,4] = 1.039845305607587
$arr2[7,5] = 1.015005973403225
$arr2[8,0] = 1.033168887651878
$arr2[8,1] = 1.031860688537889
$arr2[8,2] = 1.031881528624856
$arr2[8,3] = 1.029330908110703
$arr2[8,4] = 1.037883060730661
$arr2[8,5] = 1.01451094555118
$arr2[9,0] = 1.032949676566508
$arr2[9,1] = 1.031227782933217
$arr2[9,2] = 1.031195052730603
$arr2[9,3] = 1.028711521884123
$arr2[9,4] = 1.037031242640589
$arr2[9,5] = 1.014295690420007
$arr2[10,0] = 1.032867807140995
$arr2[10,1] = 1.030992385042313
$arr2[10,2] = 1.030939792147723
$arr2[10,3] = 1.028481209636478
$arr2[10,4] = 1.036714510094889
$arr2[10,5] = 1.014215597925435
$arr2[11,0] = 1.032885388558835
$arr2[11,1] = 1.031042892713883
$arr2[11,2] = 1.030994558821765
$arr2[11,3] = 1.028530623521708
$arr2[11,4] = 1.036782465301462
$arr2[11,5] = 1.014232784263603
$arr2[12,0] = 1.032942918294708
$arr2[12,1] = 1.031208331186986
$arr2[12,2] = 1.031173958400915
$arr2[12,3] = 1.02869248921251
$arr2[12,4] = 1.037005068183097
$arr2[12,5] = 1.014289072754416
$arr2[13,0] = 1.032978305292498
$arr2[13,1] = 1.031310222287539
$arr2[13,2] = 1.031284456124412
$arr2[13,3] = 1.028792187486487
$arr2[13,4] = 1.037142177341996
$arr2[13,5] = 1.014323735740431
$arr2[14,0] = 1.033183373699838
$arr2[14,1] = 1.031902649335702
$arr2[14,2] = 1.031927049723059
$arr2[14,3] = 1.029371980622203
$arr2[14,4] = 1.03793954722586
$arr2[14,5] = 1.014525212158393
$arr2[15,0] = 1.033311217020435
$arr2[15,1] = 1.032273717401943
$arr2[15,2] = 1.032329649761716
$arr2[15,3] = 1.029735237386167
$arr2[15,4] = 1.038439135688304
$arr2[15,5] = 1.014651349893792
$arr2[16,0] = 1.033385501435467
$arr2[16,1] = 1.032489959664549
$arr2[16,2] = 1.03256430715154
$arr2[16,3] = 1.029946964295202
$arr2[16,4] = 1.038730329844016
$arr2[16,5] = 1.014724836683386
$arr2[17,0] = 1.03341078230497
$arr2[17,1] = 1.032563659648434
$arr2[17,2] = 1.032644290135201
$arr2[17,3] = 1.030019131579274
$arr2[17,4] = 1.038829584534573
$arr2[17,5] = 1.014749879022547
$arr2[18,0] = 1.033297530074238
$arr2[18,1] = 1.032233925556487
$arr2[18,2] = 1.032286472450143
$arr2[18,3] = 1.029696279392238
$arr2[18,4] = 1.038385556064954
$arr2[18,5] = 1.014637825539053
$arr2[19,0] = 1.032925989509391
$arr2[19,1] = 1.031159622228863
$arr2[19,2] = 1.031121137261103
$arr2[19,3] = 1.028644830582059
$arr2[19,4] = 1.036939526322788
$arr2[19,5] = 1.014272500992796
$arr2[20,0] = 1.032689813688605
$arr2[20,1] = 1.030482377892504
$arr2[20,2] = 1.030386864027691
$arr2[20,3] = 1.027982325904728
$arr2[20,4] = 1.036028441386307
$arr2[20,5] = 1.014042012805997
$arr2[21,0] = 1.032815259440758
$arr2[21,1] = 1.030841568613719
$arr2[21,2] = 1.030776267341133
$arr2[21,3] = 1.028333667667285
$arr2[21,4] = 1.036511607647778
$arr2[21,5] = 1.014164274639653
$arr2[22,0] = 1.033303715494837
$arr2[22,1] = 1.032251906382359
$arr2[22,2] = 1.032305982951166
$arr2[22,3] = 1.029713883310316
$arr2[22,4] = 1.038409767030688
$arr2[22,5] = 1.014643936882239
$arr2[23,0] = 1.033859389343018
$arr2[23,1] = 1.033880870362205
$arr2[23,2] = 1.034074383574075
$arr2[23,3] = 1.031309493520469
$arr2[23,4] = 1.040604362117562
$arr2[23,5] = 1.015197141035617

$ws.Range("I2:N25").Value = $arr2

Write-Output "done"